# Apply the "Vaccine Triumph" -> "Unveiling the Wonders of the Microscopic
# Realm: A Journey into Biology" rewrite described by the commit diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Title -------------------------------------------------------------
Replace-Text "Vaccine Triumph: Science's Victory over Disease" `
             "Unveiling the Wonders of the Microscopic Realm: A Journey into Biology"

# --- Byline (collapses the three "Dr" / "." / " Sophia Cunningham" runs)
Replace-Text "Dr. Sophia Cunningham" "Sarah Williams"

# --- Email (keeps the trailing ".com" runs untouched) -------------------
Replace-Text "sophiacunningham007@gmail" "sarah.williams@ymail"

# --- Body paragraph 1 ----------------------------------------------------
Replace-Text "From the dawn of humanity, diseases have plagued humankind, leaving an indelible mark of suffering and loss" `
             "Biology, the exploration of the intricacies of life and its processes, invites us on a captivating odyssey into the enigmatic world of living organisms"

Replace-Text " The development of vaccines, however, stands as a testament to science's relentless pursuit of safeguarding human health and conquering the scourges of infectious diseases" `
             " We probe the fundamental secrets veiled within the depths of microscopic ecosystems, from cells and molecules to organisms and ecosystems"

Replace-Text " Vaccines, a beacon of hope, have transformed the landscape of public health, reducing the incidence of preventable illnesses and contributing significantly to increased life expectancy worldwide" `
             " Biology provides a lens into the marvelous tapestry of existence, revealing the wonders and complexities that shape the living realm we inhabit. As we venture deeper into this captivating subject, boundless opportunities await, orchestrating an experience of discovery and wonder"

Replace-Text "This medical revolution traces its roots to pioneering figures like Edward Jenner, who, in the 18th century, pioneered the concept of vaccination against smallpox, a deadly disease that once ravaged populations" `
             "Biology sparks our curiosity, igniting within us a burning desire to unravel the encrypted wonders concealed within the fabric of life"

Replace-Text " This breakthrough paved the way for the development of vaccines against an array of infectious agents, including polio, measles, and influenza, effectively curbing epidemics that once held humanity in their grip" `
             " We embark on a quest for knowledge and meaning that spans across multiple disciplines, weaving together insights from chemistry, physics, and mathematics to paint a comprehensive picture of the intricate mechanisms that govern the functioning of the natural world. The pursuit of biology expands our perspectives and cultivates critical-thinking skills essential for navigating the complexities of the world we inhabit"

Replace-Text "Vaccines work by introducing weakened or inactivated forms of a pathogen into the body, prompting the immune system to mount a defense" `
             "Within the spectrum of natural sciences, biology stands as a beacon of interdependence and interconnectedness, reminding us of the delicate balance and fragility of our planet's ecosystems"

Replace-Text " This process, known as immunization, equips the body with the necessary knowledge to recognize and combat the actual pathogen should it encounter it in the future, preventing illness or mitigating its severity" `
             " It unveils the ingenious mechanisms by which organisms adapt and thrive, fortifying the web of life that sustains our very existence. Understanding these intricate relationships empowers us to act as responsible stewards and ardent protectors of our natural inheritance, ensuring a sustainable, harmonious future for generations to come"

# --- Summary paragraph ----------------------------------------------------
Replace-Text "Through vaccination campaigns, entire populations have achieved herd immunity, effectively shielding vulnerable individuals and eliminating the transmission of preventable diseases" `
             "In this exploration of biology, we delved into the intricate tapestry of life, unraveling the wonders of living organisms from the molecular level to ecosystem dynamics"

Replace-Text " The eradication of smallpox, a disease that once killed millions, stands as a testament to the transformative power of vaccines" `
             " Biology's interdisciplinary nature kindles our curiosity, fostering a holistic understanding of the scientific tapestry"

Replace-Text " The near-elimination of polio and the remarkable decline in measles cases further underscore the triumph of science over disease" `
             " The insights gained from biology extend beyond the confines of scientific knowledge, fostering critical-thinking skills and shaping our perspectives on interconnectedness and sustainability"

Replace-Text " While challenges remain, such as vaccine hesitancy and the emergence of new infectious threats, the resounding success of vaccines serves as a clarion call for continued investment in research and development, ensuring the continued protection of generations to come" `
             " By embarking on this journey of biological exploration, we gain a deeper appreciation for the intricacies of the living world and our responsibility to protect the delicate balance that sustains us"

# --- Trailing empty paragraph added at the end of the document -----------
$d.Paragraphs.Add() | Out-Null
